# Refresh the crypto symbol list: update Price (col D) and Volume(1h) (col E)
# for the affected rows. Values are written as text (leading apostrophe) and
# the style is reset to "Normal" afterwards so the cells keep their original
# plain (unstyled) text representation instead of being coerced to numbers
# with an auto-applied percent/number format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.93%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'32.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'2.32%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.124"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.52%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07821"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-2.72%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.255"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-17.82%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.821"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.28%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.811"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.49%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9241"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.31%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1758"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.37%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07782"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'7.37%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08800"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.99%"
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'0.55%"
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'0.20%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001510"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.12%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005964"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.12%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.456"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-2.18%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.247"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.00%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3272"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.27%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1331"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.99%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'7.57%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.1800"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'9.08%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04605"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.22%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'1.01%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004487"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.45%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'4.66%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-1.07%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01782"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'0.71%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04793"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'6.85%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007143"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'4.31%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1364"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.51%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002126"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-3.59%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.009985"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'4.23%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006264"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-4.28%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000752"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.48%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.003605"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-58.70%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.7884"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-3.92%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002105"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.48%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002005"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.48%"
$ws.Range("E50").Style = "Normal"
